$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# A. "low speed high torque" -> "low speed, high torque"
ReplaceText "low speed high torque" "low speed, high torque"

# B. "operate the servos we only need" -> "operate the servos, we only need"
ReplaceText "operate the servos we only need" "operate the servos, we only need"

# C. "can break it, so don't put" -> "can break it…so don't put"
ReplaceText "can break it, so don" "can break it…so don"

# D. "move clockwise and counter-clockwise" -> "move clockwise or counter-clockwise"
ReplaceText "move clockwise and counter-clockwise" "move clockwise or counter-clockwise"

# E. "Connect 5 volt  pin (5V)" -> "Connect 5-volt pin (5V)"
ReplaceText "5 volt pin (5V)" "5-volt pin (5V)"

# F. "Copy/paste the code into Arduino IDE and upload to your Arduino Nano" -> with GitHub and "it"
ReplaceText "Arduino IDE and upload to your Arduino Nano" "Arduino IDE from GitHub, and upload it to your Arduino Nano"

# G. "Part 2 Power Servo using a Voltage Regulator:" -> "Part 2 Power a Servo using a Voltage Regulator:"
ReplaceText "Part 2 Power Servo using a Voltage Regulator:" "Part 2 Power a Servo using a Voltage Regulator:"

# H. "To limit the voltage we are going" -> "To limit the voltage, we are going"
ReplaceText "To limit the voltage we are going" "To limit the voltage, we are going"

# I. "with voltage regulators make sure to consult" -> "with voltage regulators, make sure to consult"
ReplaceText "with voltage regulators make sure to consult" "with voltage regulators, make sure to consult"

# J. "Another thing that to remember" -> "Another thing to remember"
ReplaceText "Another thing that to remember" "Another thing to remember"

Write-Output "Done"
